$wb = $excel.ActiveWorkbook

# --- survey sheet: insert a "phone_number" text prompt before the existing
#     send_sms row, and append two new odk_sms / odk_sms_automatic rows ---
$survey = $wb.Worksheets.Item("survey")

# Shift the old row 3 (send_sms) down to row 4, making room for the new
# phone_number prompt row that will be filled in afterwards.
$survey.Rows.Item(3).Insert()

# New row 5: odk_sms example (written first).
$survey.Range("A5").Value = "odk_sms"
$survey.Range("B5").Value = "odk_sms"
$survey.Range("C5").Value = "This will send an sms via the sms bridge."
$survey.Rows.Item(5).RowHeight = 12.75

# Row 3 (the blank row just inserted): phone_number prompt.
$survey.Range("A3").Value = "text"
$survey.Range("B3").Value = "phone_number"
$survey.Range("C3").Value = "Enter the phone number to which to send the text."
$survey.Rows.Item(3).RowHeight = 12

# New row 6: odk_sms_automatic example.
$survey.Range("A6").Value = "odk_sms_automatic"
$survey.Range("B6").Value = "odk_sms_automatic"
$survey.Range("C6").Value = "This will send an sms without requiring confirmation."
$survey.Rows.Item(6).RowHeight = 12.75

$survey.Columns.Item(2).ColumnWidth = 16.1640625

$survey.Range("C7").Select()

# --- prompt_types sheet: register the two new prompt types ---
$promptTypes = $wb.Worksheets.Item("prompt_types")

$promptTypes.Range("A3").Value = "odk_sms"
$promptTypes.Range("B3").Value = "integer"

$promptTypes.Range("A4").Value = "odk_sms_automatic"
$promptTypes.Range("B4").Value = "integer"

$promptTypes.Range("B5").Select()
$promptTypes.Activate()
